# Update the cached text of the "datetimeFigureOut" date placeholder field
# from 04/21/2023 to 04/23/2023 across the slide master and every slide
# layout (the placeholder shape is named "Date Placeholder N" on each).

$p = $ppt.ActivePresentation
$oldDate = "04/21/2023"
$newDate = "04/23/2023"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout attached to the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
